$d = $word.ActiveDocument

# The document currently ends with a single paragraph that only holds the
# "_GoBack" bookmark. Insert the three new test-case lines right before it.
$goBackPara = $d.Paragraphs($d.Paragraphs.Count)
$goBackPara.Range.InsertBefore("Display empty tree`rAdd 1 item`rRemove item from tree of length 1`r")

# Append a new, empty trailing paragraph right after the bookmark paragraph
# (which is now the second-to-last paragraph in the document).
$goBackPara = $d.Paragraphs($d.Paragraphs.Count)
$goBackPara.Range.InsertParagraphAfter()
